$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mass column (B): swap "× 10^" notation for "* 10**" notation ---
$ws.Range("B2").Value = "1.989 * 10**30"
$ws.Range("B3").Value = "3.3011*10**23"
$ws.Range("B4").Value = "4.869*10**24"
$ws.Range("B5").Value = "5.972 * 10**24"
$ws.Range("B6").Value = "6.39 * 10**23"
$ws.Range("B7").Value = "1.898 * 10**27"
$ws.Range("B8").Value = "5.683 * 10**26"
$ws.Range("B9").Value = "8.681 * 10**25"
$ws.Range("B10").Value = "1.024 * 10**26"
$ws.Range("B11").Value = "1.30900 * 10**22"

# --- Radius column (C): replace placeholder / garbled text with real numbers (km) ---
$ws.Range("C3").Value = 2439.7
$ws.Range("C4").Value = 6051.8
$ws.Range("C6").Value = 3389.5
$ws.Range("C7").Value = 69911
$ws.Range("C8").Value = 58232
$ws.Range("C9").Value = 25362
$ws.Range("C10").Value = 24622
$ws.Range("C11").Value = 1188.3

# --- Age column (D): convert text ages into real numbers (billions of years) ---
$ws.Range("D2").Value = 4.603
$ws.Range("D3").Value = 4.503
$ws.Range("D4").Value = 4.503
$ws.Range("D5").Value = 4.543
$ws.Range("D6").Value = 4.603
$ws.Range("D7").Value = 4.503
$ws.Range("D8").Value = 4.503
$ws.Range("D9").Value = 4.503
$ws.Range("D10").Value = 4.503
$ws.Range("D11").Value = 4.5

# --- New row for the tenth entry (blank placeholder body, as in source) ---
$ws.Range("B12").Value = "
"
